# Fruta / hortaliza, semanal
# Insert two new weekly-report rows (new price observations) right before the
# existing row 815 block, pushing the rest of the "Platano" data down by two
# rows (815-859 -> 817-861), and populate the two new rows with the latest
# observations (date serial 45008).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 815; everything from old row 815 downward shifts
# down to row 817 onward (Excel's default Insert shifts cells down).
$ws.Rows("815:816").Insert()

# The row that used to be 815 is now at 817 - reuse it as a style/content
# template for the columns that stay constant across this whole block
# (market/region/product metadata, units, origin, etc.).
$templateRow = 817

$cols = 1..20

foreach ($col in $cols) {
    $ws.Cells.Item(815, $col).Value = $ws.Cells.Item($templateRow, $col).Value()
    $ws.Cells.Item(816, $col).Value = $ws.Cells.Item($templateRow, $col).Value()
}

# Column D (Fecha) needs the proper date number format, same as the rest of
# the column.
$ws.Cells.Item(815, 4).NumberFormat = $ws.Cells.Item($templateRow, 4).NumberFormat()
$ws.Cells.Item(816, 4).NumberFormat = $ws.Cells.Item($templateRow, 4).NumberFormat()

# --- Row 815: "Pinton" quality ---
$ws.Cells.Item(815, 4).Value = 45008    # Fecha
$ws.Cells.Item(815, 12).Value = "Pintón"  # Calidad
$ws.Cells.Item(815, 13).Value = 200      # Volumen
$ws.Cells.Item(815, 14).Value = 21000    # Precio minimo
$ws.Cells.Item(815, 15).Value = 21000    # Precio maximo
$ws.Cells.Item(815, 16).Value = 21000    # Precio promedio ponderado
$ws.Cells.Item(815, 19).Value = 1050     # Precio $/Kg

# --- Row 816: "Primera Pinton" quality ---
$ws.Cells.Item(816, 4).Value = 45008     # Fecha
$ws.Cells.Item(816, 12).Value = "Primera Pintón"  # Calidad
$ws.Cells.Item(816, 13).Value = 300      # Volumen
$ws.Cells.Item(816, 14).Value = 22000    # Precio minimo
$ws.Cells.Item(816, 15).Value = 23000    # Precio maximo
$ws.Cells.Item(816, 16).Value = 22500    # Precio promedio ponderado
$ws.Cells.Item(816, 19).Value = 1125     # Precio $/Kg
